$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("day")

# --- Convert bsecode (column D) for rows 337-350 from text to numeric ---
$bsecodes = @{
  337 = 539448
  338 = 542726
  339 = 500790
  340 = 532978
  341 = 500043
  342 = 500570
  343 = 531642
  344 = 539336
  345 = 532955
  346 = 500877
  347 = 534816
  348 = 500547
  349 = 541153
  350 = 530965
}
foreach ($r in $bsecodes.Keys) {
  $ws.Cells.Item($r, 4).Value = $bsecodes[$r]
}

# --- Append new rows 351-365 (stock.yaml completed data) ---
# Row 351: PERSISTENT
$ws.Cells.Item(351, 1).Value = 1
$ws.Cells.Item(351, 2).Value = 'PERSISTENT'
$ws.Cells.Item(351, 3).Value = 'Persistent Systems Limited'
$ws.Cells.Item(351, 4).Value = '''533179'
$ws.Cells.Item(351, 5).Value = 1.09
$ws.Cells.Item(351, 6).Value = 4764.1
$ws.Cells.Item(351, 7).Value = 286749
$ws.Cells.Item(351, 8).Value = 'day'
$ws.Cells.Item(351, 9).Value = '14/08/2024 11:36:12'

# Row 352: JKCEMENT
$ws.Cells.Item(352, 1).Value = 2
$ws.Cells.Item(352, 2).Value = 'JKCEMENT'
$ws.Cells.Item(352, 3).Value = 'Jk Cement Limited'
$ws.Cells.Item(352, 4).Value = '''532644'
$ws.Cells.Item(352, 5).Value = -1.56
$ws.Cells.Item(352, 6).Value = 4107.2
$ws.Cells.Item(352, 7).Value = 56725
$ws.Cells.Item(352, 8).Value = 'day'
$ws.Cells.Item(352, 9).Value = '14/08/2024 11:36:12'

# Row 353: NAVINFLUOR
$ws.Cells.Item(353, 1).Value = 3
$ws.Cells.Item(353, 2).Value = 'NAVINFLUOR'
$ws.Cells.Item(353, 3).Value = 'Navin Fluorine International Limited'
$ws.Cells.Item(353, 4).Value = '''532504'
$ws.Cells.Item(353, 5).Value = -1.2
$ws.Cells.Item(353, 6).Value = 3237.6
$ws.Cells.Item(353, 7).Value = 340283
$ws.Cells.Item(353, 8).Value = 'day'
$ws.Cells.Item(353, 9).Value = '14/08/2024 11:36:12'

# Row 354: M&M
$ws.Cells.Item(354, 1).Value = 4
$ws.Cells.Item(354, 2).Value = 'M&M'
$ws.Cells.Item(354, 3).Value = 'Mahindra & Mahindra Limited'
$ws.Cells.Item(354, 4).Value = '''500520'
$ws.Cells.Item(354, 5).Value = 1
$ws.Cells.Item(354, 6).Value = 2745.25
$ws.Cells.Item(354, 7).Value = 2609359
$ws.Cells.Item(354, 8).Value = 'day'
$ws.Cells.Item(354, 9).Value = '14/08/2024 11:36:12'

# Row 355: SRF
$ws.Cells.Item(355, 1).Value = 5
$ws.Cells.Item(355, 2).Value = 'SRF'
$ws.Cells.Item(355, 3).Value = 'Srf Limited'
$ws.Cells.Item(355, 4).Value = '''503806'
$ws.Cells.Item(355, 5).Value = -1.16
$ws.Cells.Item(355, 6).Value = 2491.75
$ws.Cells.Item(355, 7).Value = 261492
$ws.Cells.Item(355, 8).Value = 'day'
$ws.Cells.Item(355, 9).Value = '14/08/2024 11:36:12'

# Row 356: MGL
$ws.Cells.Item(356, 1).Value = 6
$ws.Cells.Item(356, 2).Value = 'MGL'
$ws.Cells.Item(356, 3).Value = 'Mahanagar Gas Limited'
$ws.Cells.Item(356, 4).Value = '''539957'
$ws.Cells.Item(356, 5).Value = -3.38
$ws.Cells.Item(356, 6).Value = 1730.1
$ws.Cells.Item(356, 7).Value = 418836
$ws.Cells.Item(356, 8).Value = 'day'
$ws.Cells.Item(356, 9).Value = '14/08/2024 11:36:12'

# Row 357: OBEROIRLTY
$ws.Cells.Item(357, 1).Value = 7
$ws.Cells.Item(357, 2).Value = 'OBEROIRLTY'
$ws.Cells.Item(357, 3).Value = 'Oberoi Realty Limited'
$ws.Cells.Item(357, 4).Value = '''533273'
$ws.Cells.Item(357, 5).Value = -1.66
$ws.Cells.Item(357, 6).Value = 1716.15
$ws.Cells.Item(357, 7).Value = 651320
$ws.Cells.Item(357, 8).Value = 'day'
$ws.Cells.Item(357, 9).Value = '14/08/2024 11:36:12'

# Row 358: GNFC
$ws.Cells.Item(358, 1).Value = 8
$ws.Cells.Item(358, 2).Value = 'GNFC'
$ws.Cells.Item(358, 3).Value = 'Gujarat Narmada Valley Fertilizers And Chemicals Limited'
$ws.Cells.Item(358, 4).Value = '''500670'
$ws.Cells.Item(358, 5).Value = -1.28
$ws.Cells.Item(358, 6).Value = 649.7
$ws.Cells.Item(358, 7).Value = 3453352
$ws.Cells.Item(358, 8).Value = 'day'
$ws.Cells.Item(358, 9).Value = '14/08/2024 11:36:12'

# Row 359: BERGEPAINT
$ws.Cells.Item(359, 1).Value = 9
$ws.Cells.Item(359, 2).Value = 'BERGEPAINT'
$ws.Cells.Item(359, 3).Value = 'Berger Paints (i) Limited'
$ws.Cells.Item(359, 4).Value = '''509480'
$ws.Cells.Item(359, 5).Value = 2.31
$ws.Cells.Item(359, 6).Value = 547.2
$ws.Cells.Item(359, 7).Value = 3121387
$ws.Cells.Item(359, 8).Value = 'day'
$ws.Cells.Item(359, 9).Value = '14/08/2024 11:36:12'

# Row 360: IGL
$ws.Cells.Item(360, 1).Value = 10
$ws.Cells.Item(360, 2).Value = 'IGL'
$ws.Cells.Item(360, 3).Value = 'Indraprastha Gas Limited'
$ws.Cells.Item(360, 4).Value = '''532514'
$ws.Cells.Item(360, 5).Value = -0.28
$ws.Cells.Item(360, 6).Value = 539.15
$ws.Cells.Item(360, 7).Value = 478104
$ws.Cells.Item(360, 8).Value = 'day'
$ws.Cells.Item(360, 9).Value = '14/08/2024 11:36:12'

# Row 361: HINDPETRO
$ws.Cells.Item(361, 1).Value = 11
$ws.Cells.Item(361, 2).Value = 'HINDPETRO'
$ws.Cells.Item(361, 3).Value = 'Hindustan Petroleum Corporation Limited'
$ws.Cells.Item(361, 4).Value = '''500104'
$ws.Cells.Item(361, 5).Value = 0.47
$ws.Cells.Item(361, 6).Value = 373.1
$ws.Cells.Item(361, 7).Value = 5776570
$ws.Cells.Item(361, 8).Value = 'day'
$ws.Cells.Item(361, 9).Value = '14/08/2024 11:36:12'

# Row 362: ABFRL
$ws.Cells.Item(362, 1).Value = 12
$ws.Cells.Item(362, 2).Value = 'ABFRL'
$ws.Cells.Item(362, 3).Value = 'Aditya Birla Fashion And Retail Limited'
$ws.Cells.Item(362, 4).Value = '''535755'
$ws.Cells.Item(362, 5).Value = -0.51
$ws.Cells.Item(362, 6).Value = 311.1
$ws.Cells.Item(362, 7).Value = 1586347
$ws.Cells.Item(362, 8).Value = 'day'
$ws.Cells.Item(362, 9).Value = '14/08/2024 11:36:12'

# Row 363: ASHOKLEY
$ws.Cells.Item(363, 1).Value = 13
$ws.Cells.Item(363, 2).Value = 'ASHOKLEY'
$ws.Cells.Item(363, 3).Value = 'Ashok Leyland Limited'
$ws.Cells.Item(363, 4).Value = '''500477'
$ws.Cells.Item(363, 5).Value = -2.09
$ws.Cells.Item(363, 6).Value = 246.45
$ws.Cells.Item(363, 7).Value = 4646381
$ws.Cells.Item(363, 8).Value = 'day'
$ws.Cells.Item(363, 9).Value = '14/08/2024 11:36:12'

# Row 364: GAIL
$ws.Cells.Item(364, 1).Value = 14
$ws.Cells.Item(364, 2).Value = 'GAIL'
$ws.Cells.Item(364, 3).Value = 'Gail (india) Limited'
$ws.Cells.Item(364, 4).Value = '''532155'
$ws.Cells.Item(364, 5).Value = -0.22
$ws.Cells.Item(364, 6).Value = 226.66
$ws.Cells.Item(364, 7).Value = 12489039
$ws.Cells.Item(364, 8).Value = 'day'
$ws.Cells.Item(364, 9).Value = '14/08/2024 11:36:12'

# Row 365: GMRINFRA
$ws.Cells.Item(365, 1).Value = 15
$ws.Cells.Item(365, 2).Value = 'GMRINFRA'
$ws.Cells.Item(365, 3).Value = 'Gmr Infrastructure Limited'
$ws.Cells.Item(365, 4).Value = '''532754'
$ws.Cells.Item(365, 5).Value = -2.19
$ws.Cells.Item(365, 6).Value = 92.73
$ws.Cells.Item(365, 7).Value = 15862046
$ws.Cells.Item(365, 8).Value = 'day'
$ws.Cells.Item(365, 9).Value = '14/08/2024 11:36:12'
